$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Add new row 3 to the Logs sheet, mirroring the structure of row 2.
$logs.Range("A3").Value = "Demo inplannen"
$logs.Range("B3").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C3").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D3").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E3").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F3").Value = "2025-08-13 19:48:41"
$logs.Range("G3").Value = "Nee"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# Extend the conditional-formatting ranges on the Logs sheet so the
# existing rules also cover the newly added row 3.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $srcRange = $logs.Range("$($col)2")
    $newRange = $logs.Range("$($col)2:$($col)3")
    $fcs = $srcRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard count for the same category.
$dashboard.Range("B2").Value = 2
